{"js": "// Quick guide to the 1541 Ultimate II - work-in-progress edits:\n//  1. Move the \"_GoBack\" bookmark from the FAT16/FAT32 paragraph (end of\n//     the \"USB storage\" section intro) to the empty paragraph that sits\n//     right before the \"USB storage\" heading (this is simply where Word\n//     last left the cursor/edit point).\n//  2. Add a new sentence about checking for newer firmware into the\n//     (previously blank) first \"notes\" line at the end of the document,\n//     collapsing the now-redundant empty line that followed it.\n//  3. The document grew by one line, so the cached \"last known\" page\n//     number shown in the footer's PAGE field needs to be bumped 1 -> 2.\n\nconst body = context.document.body;\n\n// --- 1a. Remove the existing \"_GoBack\" bookmark -----------------------\n// (it currently sits right after \", or ISO files\" in the USB-storage\n// paragraph about supported file systems). No-ops harmlessly if it is\n// somehow already absent.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- 1b. Re-insert \"_GoBack\" at its new location -----------------------\n// That's the empty paragraph right before the \"USB storage\" heading.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet blankBeforeUsbStorage = null;\nfor (let i = 0; i < paragraphs.items.length - 1; i++) {\n  const current = paragraphs.items[i];\n  const next = paragraphs.items[i + 1];\n  if (current.text === \"\" && next.text === \"USB storage\") {\n    blankBeforeUsbStorage = current;\n    break;\n  }\n}\n\nif (!blankBeforeUsbStorage) {\n  throw new Error(\"Could not locate the blank paragraph before 'USB storage'.\");\n}\n\nblankBeforeUsbStorage.getRange().insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- 2. Fill in the first blank \"notes\" line, drop the next blank line -\n// The \"Notes\" section near the end has several ruled lines made from a\n// lone tab character on a paragraph with an underscored right tab stop.\n// The very first of those gets the new firmware note; the blank line\n// that used to follow it is removed so the layout doesn't gain an extra\n// row.\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet notesHeadingIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"Notes\") {\n    notesHeadingIndex = i;\n    break;\n  }\n}\n\nif (notesHeadingIndex === -1) {\n  throw new Error(\"Could not locate the 'Notes' heading paragraph.\");\n}\n\n// First paragraph after the heading that contains only a tab character.\nlet tabLineIndex = -1;\nfor (let i = notesHeadingIndex + 1; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"\\t\") {\n    tabLineIndex = i;\n    break;\n  }\n}\n\nif (tabLineIndex === -1) {\n  throw new Error(\"Could not locate the blank ruled note line.\");\n}\n\nconst noteParagraph = paragraphs.items[tabLineIndex];\nnoteParagraph.insertText(\n  \"At the time of arrival of this 1541 Ultimate-II unit, there might be a newer version of the firmware available for download from the website, with enhanced functionality.\",\n  Word.InsertLocation.replace\n);\n\nconst blankAfterNote = paragraphs.items[tabLineIndex + 1];\nblankAfterNote.delete();\nawait context.sync();\n\n// --- 3. Update the cached PAGE field result in the primary footer -----\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nconst primaryFooter = sections.items[0].getFooter(Word.HeaderFooterType.primary);\nconst pageNumberHits = primaryFooter.search(\"1\", { matchWholeWord: true });\npageNumberHits.load(\"items\");\nawait context.sync();\n\nif (pageNumberHits.items.length === 0) {\n  throw new Error(\"Could not locate the cached page-number field result.\");\n}\n\npageNumberHits.items[0].insertText(\"2\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Quick guide to the 1541 Ultimate II - work-in-progress edits:\n#  1. Move the \"_GoBack\" bookmark from the FAT16/FAT32 paragraph (end of\n#     the \"USB storage\" section intro) to the empty paragraph that sits\n#     right before the \"USB storage\" heading (this is simply where Word\n#     last left the cursor/edit point).\n#  2. Add a new sentence about checking for newer firmware into the\n#     (previously blank) first \"notes\" line at the end of the document,\n#     collapsing the now-redundant empty line that followed it.\n#  3. The document grew by one line, so the cached \"last known\" page\n#     number shown in the footer's PAGE field needs to be bumped 1 -> 2.\n\n$d = $word.ActiveDocument\n\n# --- 1a. Remove the existing \"_GoBack\" bookmark ------------------------\n$existing = $d.Bookmarks(\"_GoBack\")\nif ($existing -ne $null) {\n    $existing.Delete()\n}\n\n# --- 1b. Re-insert \"_GoBack\" at its new location ------------------------\n# That's the empty paragraph right before the \"USB storage\" heading.\n# NOTE: Paragraph.Range.Text always carries a trailing paragraph-mark\n# (\"`r\"), so compare against that rather than a bare empty string.\n$target = $null\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -lt $count; $i++) {\n    $current = $d.Paragraphs($i)\n    $next = $d.Paragraphs($i + 1)\n    if ($current.Range.Text -eq \"`r\" -and $next.Range.Text -eq \"USB storage`r\") {\n        $target = $current\n        break\n    }\n}\nif ($target -eq $null) {\n    throw \"Could not locate the blank paragraph before 'USB storage'.\"\n}\n$d.Bookmarks.Add(\"_GoBack\", $target.Range)\n\n# --- 2. Fill in the first blank \"notes\" line, drop the next blank line -\n# The \"Notes\" section near the end has several ruled lines made from a\n# lone tab character on a paragraph with an underscored right tab stop.\n# The very first of those gets the new firmware note; the blank line\n# that used to follow it is removed so the layout doesn't gain an extra\n# row.\n$notesIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    if ($d.Paragraphs($i).Range.Text -eq \"Notes`r\") {\n        $notesIndex = $i\n        break\n    }\n}\nif ($notesIndex -eq -1) {\n    throw \"Could not locate the 'Notes' heading paragraph.\"\n}\n\n$tabLineIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = $notesIndex + 1; $i -le $count; $i++) {\n    if ($d.Paragraphs($i).Range.Text -eq \"`t`r\") {\n        $tabLineIndex = $i\n        break\n    }\n}\nif ($tabLineIndex -eq -1) {\n    throw \"Could not locate the blank ruled note line.\"\n}\n\n$noteParagraph = $d.Paragraphs($tabLineIndex)\n$noteParagraph.Range.Text = \"At the time of arrival of this 1541 Ultimate-II unit, there might be a newer version of the firmware available for download from the website, with enhanced functionality.\"\n\n$blankAfterNote = $d.Paragraphs($tabLineIndex + 1)\n$blankAfterNote.Range.Delete()\n\n# --- 3. Update the cached PAGE field result in the primary footer -----\n$section = $d.Sections(1)\n$footer = $section.Footers(1)\n$firstChar = $footer.Range.Characters.Item(1)\n$firstChar.Text = \"2\"\n"}
